$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update header labels (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Update column C values (rows 2-25) ---
$cValues = @{
    2  = 1909.084588129339
    3  = 4729.735976516416
    4  = 471.181692645893
    5  = 5082.354756663512
    6  = 1955.461557360978
    7  = 492.3430015592067
    8  = 2024.117324382548
    9  = 5360.226632400601
    10 = 513.7390871590731
    11 = 2094.024217383061
    12 = 5642.578115155247
    13 = 534.5063430177229
    14 = 2201.396847776877
    15 = 5919.20956823756
    16 = 2286.013198234259
    17 = 558.2093442539386
    18 = 5996.49696468919
    19 = 2361.056581219794
    20 = 579.0880693780265
    21 = 6114.227214287786
    22 = 2425.561644739583
    23 = 584.2111078769213
    24 = 6262.368904654469
    25 = 1431.756130822538
}

foreach ($row in $cValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
}
